$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values per repulled data / mean calculation fix
$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = -4
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 1
